$d = $word.ActiveDocument

# Locate the "Michael Thompson" text in the document
$rng = $d.Content
$found = $rng.Find.Execute("Michael Thompson", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse the range to its end point (right after "Michael Thompson")
    $rng.Collapse(0)

    # Insert the new text as a new run right after it
    $rng.InsertAfter(" 222333")

    # The inserted text now occupies $rng (InsertAfter updates the range to span the inserted text)
    $rng.Font.Name = "Aptos"
    $rng.Font.NameAscii = "Aptos"
    $rng.Font.NameFarEast = "Aptos"
    $rng.Font.NameOther = "Aptos"
    $rng.Font.NameComplexScript = "Segoe UI"
    $rng.Font.Size = 14
    $rng.LanguageID = 1033
}
